$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.1
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 2.9
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 6.1
